# Update market-price / profit figures (columns H-N) on a handful of
# leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets,
# reflecting refreshed marketboard data from the scheduled scraper run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 305.5
$ws.Range("I2").Value = 74.833336
$ws.Range("J2").Value = 997.5
$ws.Range("K2").Value = 74.833336
$ws.Range("L2").Value = 997.5
$ws.Range("M2").Value = 38.166664
$ws.Range("N2").Value = -1223.5

$ws.Range("H28").Value = 1699.3334
$ws.Range("J28").Value = 3112.4285
$ws.Range("L28").Value = 3112.4285
$ws.Range("N28").Value = -4082.4285

$ws.Range("H64").Value = 5000
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4752

$ws.Range("H67").Value = 5000
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4142

$ws.Range("H70").Value = 51960.277
$ws.Range("I70").Value = 128714.71
$ws.Range("J70").Value = 3116.5454
$ws.Range("K70").Value = 386144.13
$ws.Range("L70").Value = 9349.636200000001
$ws.Range("M70").Value = -385874.13
$ws.Range("N70").Value = -9889.636200000001

$ws.Range("H73").Value = 51960.277
$ws.Range("I73").Value = 128714.71
$ws.Range("J73").Value = 3116.5454
$ws.Range("K73").Value = 386144.13
$ws.Range("L73").Value = 9349.636200000001
$ws.Range("M73").Value = -385208.13
$ws.Range("N73").Value = -11221.6362

$ws.Range("H86").Value = 1224.75
$ws.Range("I86").Value = 633
$ws.Range("K86").Value = 633
$ws.Range("M86").Value = 490

$ws.Range("H89").Value = 1224.75
$ws.Range("I89").Value = 633
$ws.Range("K89").Value = 3165
$ws.Range("M89").Value = 2451

$ws.Range("H101").Value = 8494.333000000001
$ws.Range("I101").Value = 5773.2
$ws.Range("K101").Value = 17319.6
$ws.Range("M101").Value = -15697.6

$ws.Range("H107").Value = 1336.1666
$ws.Range("I107").Value = 270.2
$ws.Range("K107").Value = 270.2
$ws.Range("M107").Value = 1649.8

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H125").Value = 1917
$ws.Range("I125").Value = 2151.5
$ws.Range("K125").Value = 19363.5
$ws.Range("M125").Value = -16903.5

$ws.Range("H138").Value = 2874.7808
$ws.Range("I138").Value = 4237.778
$ws.Range("K138").Value = 12713.334
$ws.Range("M138").Value = -7573.334000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 979.5714
$ws.Range("J4").Value = 979.5
$ws.Range("L4").Value = 979.5
$ws.Range("N4").Value = -1211.5

$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H74").Value = 3228.0715
$ws.Range("I74").Value = 2855
$ws.Range("J74").Value = 4596
$ws.Range("K74").Value = 2855
$ws.Range("L74").Value = 4596
$ws.Range("M74").Value = -1981
$ws.Range("N74").Value = -6344

$ws.Range("H77").Value = 3228.0715
$ws.Range("I77").Value = 2855
$ws.Range("J77").Value = 4596
$ws.Range("K77").Value = 14275
$ws.Range("L77").Value = 22980
$ws.Range("M77").Value = -9907
$ws.Range("N77").Value = -31716

$ws.Range("H132").Value = 1446.8667
$ws.Range("I132").Value = 1273.2273
$ws.Range("K132").Value = 3819.6819
$ws.Range("M132").Value = -1289.6819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1556.8182
$ws.Range("J134").Value = 2245.5
$ws.Range("L134").Value = 6736.5
$ws.Range("N134").Value = -11806.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3114.4707
$ws.Range("I58").Value = 3089.9167
$ws.Range("K58").Value = 3089.9167
$ws.Range("M58").Value = -2886.9167

$ws.Range("H63").Value = 70271
$ws.Range("J63").Value = 70271
$ws.Range("L63").Value = 70271
$ws.Range("N63").Value = -71643

$ws.Range("H66").Value = 70271
$ws.Range("J66").Value = 70271
$ws.Range("L66").Value = 210813
$ws.Range("N66").Value = -217677

$ws.Range("H136").Value = 3114.4707
$ws.Range("I136").Value = 3089.9167
$ws.Range("K136").Value = 9269.750100000001
$ws.Range("M136").Value = -6719.750100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 591.13635
$ws.Range("I113").Value = 743.5
$ws.Range("K113").Value = 2230.5
$ws.Range("M113").Value = -60.5

$ws.Range("H134").Value = 619.1429000000001
$ws.Range("I134").Value = 619.1429000000001
$ws.Range("K134").Value = 1857.4287
$ws.Range("M134").Value = 3212.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 348.5
$ws.Range("I2").Value = 415.2
$ws.Range("K2").Value = 415.2
$ws.Range("M2").Value = -302.2

$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -860
$ws.Range("N12").ClearContents()

$ws.Range("H14").Value = 25014600
$ws.Range("I14").Value = 31250000
$ws.Range("K14").Value = 31250000
$ws.Range("M14").Value = -31249832

$ws.Range("H80").Value = 12780.529
$ws.Range("I80").Value = 7198
$ws.Range("J80").Value = 17742.777
$ws.Range("K80").Value = 7198
$ws.Range("L80").Value = 17742.777
$ws.Range("M80").Value = -6200
$ws.Range("N80").Value = -19738.777

$ws.Range("H83").Value = 12780.529
$ws.Range("I83").Value = 7198
$ws.Range("J83").Value = 17742.777
$ws.Range("K83").Value = 35990
$ws.Range("L83").Value = 88713.88499999999
$ws.Range("M83").Value = -30998
$ws.Range("N83").Value = -98697.88499999999

$ws.Range("H107").Value = 876.3077
$ws.Range("I107").Value = 470.8
$ws.Range("J107").Value = 1129.75
$ws.Range("K107").Value = 470.8
$ws.Range("L107").Value = 1129.75
$ws.Range("M107").Value = 1449.2
$ws.Range("N107").Value = -4969.75

$ws.Range("H132").Value = 2999.6667
$ws.Range("I132").Value = 2599
$ws.Range("K132").Value = 7797
$ws.Range("M132").Value = -5267

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2190
$ws.Range("I7").Value = 1432.25
$ws.Range("J7").Value = 3402.4
$ws.Range("K7").Value = 1432.25
$ws.Range("L7").Value = 3402.4
$ws.Range("M7").Value = -1320.25
$ws.Range("N7").Value = -3626.4

$ws.Range("H16").Value = 392.72726
$ws.Range("I16").Value = 392.72726
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 392.72726
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -222.72726
$ws.Range("N16").ClearContents()

$ws.Range("H82").Value = 1787.3077
$ws.Range("I82").Value = 1812.3636
$ws.Range("J82").Value = 1649.5
$ws.Range("K82").Value = 1812.3636
$ws.Range("L82").Value = 1649.5
$ws.Range("M82").Value = -1451.3636
$ws.Range("N82").Value = -2371.5

$ws.Range("H85").Value = 1787.3077
$ws.Range("I85").Value = 1812.3636
$ws.Range("J85").Value = 1649.5
$ws.Range("K85").Value = 1812.3636
$ws.Range("L85").Value = 1649.5
$ws.Range("M85").Value = -564.3635999999999
$ws.Range("N85").Value = -4145.5

$ws.Range("H126").Value = 2190
$ws.Range("I126").Value = 1432.25
$ws.Range("J126").Value = 3402.4
$ws.Range("K126").Value = 4296.75
$ws.Range("L126").Value = 10207.2
$ws.Range("M126").Value = -1826.75
$ws.Range("N126").Value = -15147.2

$ws.Range("H136").Value = 2468.88
$ws.Range("I136").Value = 2496.5454
$ws.Range("J136").Value = 2266
$ws.Range("K136").Value = 7489.6362
$ws.Range("L136").Value = 6798
$ws.Range("M136").Value = -4939.6362
$ws.Range("N136").Value = -11898

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14720.4
$ws.Range("I62").Value = 14720.4
$ws.Range("K62").Value = 14720.4
$ws.Range("M62").Value = -14096.4

$ws.Range("H65").Value = 14720.4
$ws.Range("I65").Value = 14720.4
$ws.Range("K65").Value = 73602
$ws.Range("M65").Value = -70482

$ws.Range("H109").Value = 93100
$ws.Range("J109").Value = 93100
$ws.Range("L109").Value = 93100
$ws.Range("N109").Value = -95874

$ws.Range("H113").Value = 471.625
$ws.Range("I113").Value = 418.75
$ws.Range("J113").Value = 524.5
$ws.Range("K113").Value = 1256.25
$ws.Range("L113").Value = 1573.5
$ws.Range("M113").Value = 913.75
$ws.Range("N113").Value = -5913.5
